$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.910.62"
$ws.Range("E2").Value = "  -1.65%  "

# Row 3
$ws.Range("D3").Value = "2.285.63"
$ws.Range("E3").Value = "  -2.63%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.21"
$ws.Range("E5").Value = "  -3.62%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.80"
$ws.Range("E6").Value = "  +3.76%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.619"
$ws.Range("E7").Value = "  -2.77%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.605"
$ws.Range("E9").Value = "  -2.32%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.34"
$ws.Range("E10").Value = "  +1.04%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0909"
$ws.Range("E11").Value = "  -1.13%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.25"
$ws.Range("E12").Value = "  -1.73%  "

# Row 13
$ws.Range("E13").Value = "  +0.30%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.969"
$ws.Range("E14").Value = "  -2.50%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.46"
$ws.Range("E15").Value = "  -4.29%  "

# Row 16
$ws.Range("D16").Value = "2.632.70"
$ws.Range("E16").Value = "  -2.57%  "

# Row 17
$ws.Range("D17").Value = "2.288.86"
$ws.Range("E17").Value = "  -2.52%  "

# Row 18
$ws.Range("D18").Value = "42.017.09"
$ws.Range("E18").Value = "  -1.23%  "

# Row 19
$ws.Range("E19").Value = "  -4.64%  "

# Row 20
$ws.Range("E20").Value = "  -1.46%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "74.27"
$ws.Range("E21").Value = "  -2.26%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.45"
$ws.Range("E22").Value = "  -6.48%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "256.92"
$ws.Range("E23").Value = "  -2.93%  "

# Row 24
$ws.Range("E24").Value = "  +0.05%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.31"
$ws.Range("E25").Value = "  -6.63%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.29%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.94"
$ws.Range("E27").Value = "  -4.16%  "

# Row 28
$ws.Range("E28").Value = "  +3.41%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.72"
$ws.Range("E29").Value = "  -0.19%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.15"

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.75"
$ws.Range("E31").Value = "  +2.06%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0891"
$ws.Range("E32").Value = "  -0.57%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.90"
$ws.Range("E33").Value = "  -6.17%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.79"
$ws.Range("E34").Value = "  -3.93%  "

# Row 35 and 36 swap: Kaspa <-> Stellar
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.129"
$ws.Range("E35").Value = "  -1.93%  "

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.118"
$ws.Range("E36").Value = "  +10.49%  "

# Row 37
$ws.Range("E37").Value = "  +0.29%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0351"
$ws.Range("E38").Value = "  -1.21%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.75"
$ws.Range("E39").Value = "  -2.94%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.61"
$ws.Range("E40").Value = "  -3.99%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.08"
$ws.Range("E41").Value = "  +3.60%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.32"
$ws.Range("E42").Value = "  +7.37%  "

# Row 43
$ws.Range("E43").Value = "  -2.72%  "

# Row 44
$ws.Range("E44").Value = "  -3.69%  "

# Row 45
$ws.Range("E45").Value = "  +0.21%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.27"
$ws.Range("E46").Value = "  +3.59%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "111.84"
$ws.Range("E47").Value = "  -6.92%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.02"
$ws.Range("E48").Value = "  -1.35%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.31"
$ws.Range("E49").Value = "  -3.55%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.18"
$ws.Range("E50").Value = "  +4.44%  "

# Row 51
$ws.Range("D51").Value = "1.561.08"
$ws.Range("E51").Value = "  +0.71%  "
